$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.754.00"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "3.172.25"
$ws.Range("E3").Value = "  -4.74%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.15"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.33"
$ws.Range("E6").Value = "  -2.70%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").Value = "3.170.04"
$ws.Range("E9").Value = "  -4.70%  "

$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("E11").Value = "  -3.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.393"
$ws.Range("E12").Value = "  -3.20%  "

$ws.Range("D13").Value = "3.720.88"
$ws.Range("E13").Value = "  -4.67%  "

$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.23"
$ws.Range("E15").Value = "  -4.06%  "

$ws.Range("D16").Value = "65.700.63"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("E17").Value = "  -2.11%  "

$ws.Range("D18").Value = "3.172.65"
$ws.Range("E18").Value = "  -4.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.94"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.496"
$ws.Range("E25").Value = "  -3.91%  "

$ws.Range("D26").Value = "3.305.80"
$ws.Range("E26").Value = "  -4.85%  "

$ws.Range("E27").Value = "  -6.38%  "

$ws.Range("E28").Value = "  +3.78%  "

$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("E32").Value = "  -0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.38"
$ws.Range("E33").Value = "  -3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.09"
$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.63"
$ws.Range("E35").Value = "  -2.33%  "

$ws.Range("E36").Value = "  -0.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.89"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.836"
$ws.Range("E39").Value = "  -0.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.35"
$ws.Range("E41").Value = "  -3.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("D43").Value = "2.651.24"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.19"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.21"
$ws.Range("E45").Value = "  -1.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.81"
$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0658"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "329.31"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.83"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0274"
$ws.Range("E50").Value = "  -1.06%  "

$ws.Range("E51").Value = "  -0.63%  "
